$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1806.4387
$ws.Range("I15").Value = 1806.4387
$ws.Range("K15").Value = 5419.3161
$ws.Range("M15").Value = -5250.3161
$ws.Range("H69").Value = 5274.7144
$ws.Range("J69").Value = 4979.4
$ws.Range("L69").Value = 14938.2
$ws.Range("N69").Value = -16686.2
$ws.Range("H72").Value = 5274.7144
$ws.Range("J72").Value = 4979.4
$ws.Range("L72").Value = 44814.6
$ws.Range("N72").Value = -53550.6
$ws.Range("H94").Value = 1684.1666
$ws.Range("I94").Value = 1684.1666
$ws.Range("K94").Value = 1684.1666
$ws.Range("M94").Value = -1233.1666
$ws.Range("H100").Value = 919.0454999999999
$ws.Range("I100").Value = 853.6
$ws.Range("J100").Value = 1059.2858
$ws.Range("K100").Value = 853.6
$ws.Range("L100").Value = 1059.2858
$ws.Range("M100").Value = -312.6
$ws.Range("N100").Value = -2141.2858
$ws.Range("H127").Value = 736.875
$ws.Range("I127").Value = 632.5
$ws.Range("J127").Value = 1050
$ws.Range("K127").Value = 1897.5
$ws.Range("L127").Value = 3150
$ws.Range("M127").Value = 3062.5
$ws.Range("N127").Value = -13070
$ws.Range("H129").Value = 918.2833000000001
$ws.Range("J129").Value = 973.63635
$ws.Range("L129").Value = 2920.90905
$ws.Range("N129").Value = -12920.90905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1779.0952
$ws.Range("I45").Value = 1721.875
$ws.Range("J45").Value = 1962.2
$ws.Range("K45").Value = 1721.875
$ws.Range("L45").Value = 1962.2
$ws.Range("M45").Value = -1344.875
$ws.Range("N45").Value = -2716.2
$ws.Range("H74").Value = 13542.474
$ws.Range("I74").Value = 9832.866
$ws.Range("J74").Value = 27453.5
$ws.Range("K74").Value = 9832.866
$ws.Range("L74").Value = 27453.5
$ws.Range("M74").Value = -8958.866
$ws.Range("N74").Value = -29201.5
$ws.Range("H77").Value = 13542.474
$ws.Range("I77").Value = 9832.866
$ws.Range("J77").Value = 27453.5
$ws.Range("K77").Value = 49164.33
$ws.Range("L77").Value = 137267.5
$ws.Range("M77").Value = -44796.33
$ws.Range("N77").Value = -146003.5
$ws.Range("H97").Value = 1475.2354
$ws.Range("I97").Value = 1210.8182
$ws.Range("J97").Value = 1960
$ws.Range("K97").Value = 1210.8182
$ws.Range("L97").Value = 1960
$ws.Range("M97").Value = -714.8181999999999
$ws.Range("N97").Value = -2952
$ws.Range("H102").Value = 2170.8096
$ws.Range("I102").Value = 1646.3529
$ws.Range("J102").Value = 4399.75
$ws.Range("K102").Value = 1646.3529
$ws.Range("L102").Value = 4399.75
$ws.Range("M102").Value = -24.35290000000009
$ws.Range("N102").Value = -7643.75
$ws.Range("H122").Value = 3379047
$ws.Range("I122").Value = 687.1667
$ws.Range("J122").Value = 125000000
$ws.Range("K122").Value = 2061.5001
$ws.Range("L122").Value = 375000000
$ws.Range("M122").Value = 388.4998999999998
$ws.Range("N122").Value = -375004900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2285.5715
$ws.Range("I94").Value = 2249.75
$ws.Range("J94").Value = 2333.3333
$ws.Range("K94").Value = 2249.75
$ws.Range("L94").Value = 2333.3333
$ws.Range("M94").Value = -1798.75
$ws.Range("N94").Value = -3235.3333
$ws.Range("H99").Value = 1632
$ws.Range("I99").Value = 1574.5454
$ws.Range("J99").Value = 1702.2222
$ws.Range("K99").Value = 1574.5454
$ws.Range("L99").Value = 1702.2222
$ws.Range("M99").Value = -76.54539999999997
$ws.Range("N99").Value = -4698.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3519.951
$ws.Range("I31").Value = 3707.6191
$ws.Range("J31").Value = 3105.1052
$ws.Range("K31").Value = 3707.6191
$ws.Range("L31").Value = 3105.1052
$ws.Range("M31").Value = -3412.6191
$ws.Range("N31").Value = -3695.1052
$ws.Range("H34").Value = 3519.951
$ws.Range("I34").Value = 3707.6191
$ws.Range("J34").Value = 3105.1052
$ws.Range("K34").Value = 3707.6191
$ws.Range("L34").Value = 3105.1052
$ws.Range("M34").Value = -3505.6191
$ws.Range("N34").Value = -3509.1052
$ws.Range("H105").Value = 780.875
$ws.Range("I105").Value = 554.5
$ws.Range("J105").Value = 2365.5
$ws.Range("K105").Value = 554.5
$ws.Range("L105").Value = 2365.5
$ws.Range("M105").Value = 1192.5
$ws.Range("N105").Value = -5859.5
$ws.Range("H107").Value = 1112.381
$ws.Range("I107").Value = 1418.4
$ws.Range("J107").Value = 834.1818
$ws.Range("K107").Value = 1418.4
$ws.Range("L107").Value = 834.1818
$ws.Range("M107").Value = 501.5999999999999
$ws.Range("N107").Value = -4674.1818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 633.9231
$ws.Range("J92").Value = 649.875
$ws.Range("L92").Value = 1949.625
$ws.Range("N92").Value = -4445.625
$ws.Range("H131").Value = 30103.303
$ws.Range("I131").Value = 1342.1428
$ws.Range("J131").Value = 51295.74
$ws.Range("K131").Value = 4026.4284
$ws.Range("L131").Value = 153887.22
$ws.Range("M131").Value = 1013.5716
$ws.Range("N131").Value = -163967.22

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 14018
$ws.Range("J40").Value = 14018
$ws.Range("L40").Value = 14018
$ws.Range("N40").Value = -14320
$ws.Range("H70").Value = 5879.6665
$ws.Range("I70").Value = 5372.7036
$ws.Range("J70").Value = 6450
$ws.Range("K70").Value = 5372.7036
$ws.Range("L70").Value = 6450
$ws.Range("M70").Value = -5102.7036
$ws.Range("N70").Value = -6990
$ws.Range("H73").Value = 5879.6665
$ws.Range("I73").Value = 5372.7036
$ws.Range("J73").Value = 6450
$ws.Range("K73").Value = 5372.7036
$ws.Range("L73").Value = 6450
$ws.Range("M73").Value = -4436.7036
$ws.Range("N73").Value = -8322
$ws.Range("H97").Value = 1325.2727
$ws.Range("I97").Value = 1099.6666
$ws.Range("J97").Value = 1596
$ws.Range("K97").Value = 1099.6666
$ws.Range("L97").Value = 1596
$ws.Range("M97").Value = -603.6666
$ws.Range("N97").Value = -2588
$ws.Range("H132").Value = 3764.7144
$ws.Range("I132").Value = 1691.4103
$ws.Range("J132").Value = 11850.6
$ws.Range("K132").Value = 5074.2309
$ws.Range("L132").Value = 35551.8
$ws.Range("M132").Value = -2544.2309
$ws.Range("N132").Value = -40611.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 797.5454999999999
$ws.Range("I22").Value = 658.875
$ws.Range("K22").Value = 658.875
$ws.Range("M22").Value = -363.875
$ws.Range("H27").Value = 797.5454999999999
$ws.Range("I27").Value = 658.875
$ws.Range("K27").Value = 658.875
$ws.Range("M27").Value = -551.875
$ws.Range("H46").Value = 997.7778
$ws.Range("I46").Value = 376
$ws.Range("J46").Value = 1775
$ws.Range("K46").Value = 376
$ws.Range("L46").Value = 1775
$ws.Range("M46").Value = -188
$ws.Range("N46").Value = -2151
$ws.Range("H82").Value = 2000
$ws.Range("J82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2000
$ws.Range("J85").Value = 3000
$ws.Range("L85").Value = 3000
$ws.Range("N85").Value = -5496
$ws.Range("H93").Value = 465.61905
$ws.Range("I93").Value = 438.14285
$ws.Range("J93").Value = 520.5714
$ws.Range("K93").Value = 438.14285
$ws.Range("L93").Value = 520.5714
$ws.Range("M93").Value = 809.85715
$ws.Range("N93").Value = -3016.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3400
$ws.Range("I81").Value = 1940
$ws.Range("J81").Value = 4312.5
$ws.Range("K81").Value = 3880
$ws.Range("L81").Value = 8625
$ws.Range("M81").Value = -2819
$ws.Range("N81").Value = -10747
$ws.Range("H84").Value = 3400
$ws.Range("I84").Value = 1940
$ws.Range("J84").Value = 4312.5
$ws.Range("K84").Value = 19400
$ws.Range("L84").Value = 43125
$ws.Range("M84").Value = -14096
$ws.Range("N84").Value = -53733
